# Auto-generated PowerShell COM-interop script reproducing the
# QuestionBank-Template-V2.xlsx edit: new MultipleChoice + Matching rows.

$wb = $excel.ActiveWorkbook

# --- MultipleChoice: append question rows 19-36 ---
$wsMC = $wb.Worksheets.Item("MultipleChoice")

# Row 19
$wsMC.Range("A19").Value = "Wich of the follwing is a fruit?"
$wsMC.Range("B19").Value = "Applesss"
$wsMC.Range("C19").Value = "Banana"
$wsMC.Range("D19").Value = "Carrot"
$wsMC.Range("E19").Value = "Dog"
$wsMC.Range("F19").Value = "None of the aboves"
$wsMC.Range("G19").Value = "a"
$wsMC.Range("H19").Value = "fruit.png"
$wsMC.Range("I19").Value = "food"

# Row 20
$wsMC.Range("A20").Value = "What animal do bark?"
$wsMC.Range("B20").Value = "Cat"
$wsMC.Range("C20").Value = "Snake"
$wsMC.Range("D20").Value = "Dog"
$wsMC.Range("E20").Value = "Fish"
$wsMC.Range("F20").Value = "All of the Above"
$wsMC.Range("G20").Value = "c"
$wsMC.Range("H20").Value = "dog.png"
$wsMC.Range("I20").Value = "animal"

# Row 21
$wsMC.Range("A21").Value = "Where is Eiffel Towerr?"
$wsMC.Range("B21").Value = "Paris"
$wsMC.Range("C21").Value = "London"
$wsMC.Range("D21").Value = "Bkk"
$wsMC.Range("E21").Value = "Your Home"
$wsMC.Range("F21").Value = "Russia"
$wsMC.Range("G21").Value = "a"
$wsMC.Range("H21").Value = "tower.png"
$wsMC.Range("I21").Value = "places"

# Row 22
$wsMC.Range("A22").Value = "Which is the worssst fruit ever?"
$wsMC.Range("B22").Value = "Grape"
$wsMC.Range("C22").Value = "Mango"
$wsMC.Range("D22").Value = "Appless"
$wsMC.Range("E22").Value = "Durian"
$wsMC.Range("F22").Value = "Lettuce"
$wsMC.Range("G22").Value = "d"
$wsMC.Range("H22").Value = "durian.png"
$wsMC.Range("I22").Value = "food"

# Row 23
$wsMC.Range("A23").Value = "Which of this is a public park?"
$wsMC.Range("B23").Value = "Greenfeild"
$wsMC.Range("C23").Value = "School"
$wsMC.Range("D23").Value = "Metro Plaza"
$wsMC.Range("E23").Value = "City Librry"
$wsMC.Range("F23").Value = "Cinema"
$wsMC.Range("G23").Value = "a"
$wsMC.Range("H23").Value = "park.png"
$wsMC.Range("I23").Value = "places"

# Row 24
$wsMC.Range("A24").Value = "Which of this is a public park?"
$wsMC.Range("B24").Value = "Parksss"
$wsMC.Range("C24").Value = "Home"
$wsMC.Range("D24").Value = "Forest"
$wsMC.Range("E24").Value = "Train Stashion"
$wsMC.Range("F24").Value = "Theaterss"
$wsMC.Range("G24").Value = "a"
$wsMC.Range("H24").Value = "park2.png"
$wsMC.Range("I24").Value = "places"

# Row 25
$wsMC.Range("A25").Value = "How many corners has a squar?"
$wsMC.Range("B25").Value = 1
$wsMC.Range("C25").Value = 2
$wsMC.Range("D25").Value = 3
$wsMC.Range("E25").Value = 4
$wsMC.Range("F25").Value = 5
$wsMC.Range("G25").Value = "d"
$wsMC.Range("H25").Value = "square2.png"
$wsMC.Range("I25").Value = "maths"

# Row 26
$wsMC.Range("A26").Value = "How many corner are in a square?"
$wsMC.Range("B26").Value = "One"
$wsMC.Range("C26").Value = "Two"
$wsMC.Range("D26").Value = "Three"
$wsMC.Range("E26").Value = "Four"
$wsMC.Range("F26").Value = "None"
$wsMC.Range("G26").Value = "d"
$wsMC.Range("H26").Value = "square3.png"
$wsMC.Range("I26").Value = "maths"

# Row 27
$wsMC.Range("A27").Value = "Which of the follwing is an animall?"
$wsMC.Range("B27").Value = "Dog"
$wsMC.Range("C27").Value = "Cat"
$wsMC.Range("D27").Value = "Tree"
$wsMC.Range("E27").Value = "Rice"
$wsMC.Range("F27").Value = "Apple"
$wsMC.Range("G27").Value = "a"
$wsMC.Range("H27").Value = "anim1.png"
$wsMC.Range("I27").Value = "animal"

# Row 28
$wsMC.Range("A28").Value = "Wich is your favorit place?"
$wsMC.Range("B28").Value = "Parkk"
$wsMC.Range("C28").Value = "Homess"
$wsMC.Range("D28").Value = "Malll"
$wsMC.Range("E28").Value = "Pooll"
$wsMC.Range("F28").Value = "Theatersss"
$wsMC.Range("G28").Value = "a"
$wsMC.Range("H28").Value = "place1.png"
$wsMC.Range("I28").Value = "places"

# Row 29
$wsMC.Range("A29").Value = "What animal can flying?"
$wsMC.Range("B29").Value = "Dog"
$wsMC.Range("C29").Value = "Cat"
$wsMC.Range("D29").Value = "Bird"
$wsMC.Range("E29").Value = "Fish"
$wsMC.Range("F29").Value = "Horse"
$wsMC.Range("G29").Value = "c"
$wsMC.Range("H29").Value = "bird2.png"
$wsMC.Range("I29").Value = "animal"

# Row 30
$wsMC.Range("A30").Value = "What is the badest fruit?"
$wsMC.Range("B30").Value = "Banana"
$wsMC.Range("C30").Value = "Applle"
$wsMC.Range("D30").Value = "Mango"
$wsMC.Range("E30").Value = "Durian"
$wsMC.Range("F30").Value = "Potato"
$wsMC.Range("G30").Value = "d"
$wsMC.Range("H30").Value = "fruit2.png"
$wsMC.Range("I30").Value = "food"

# Row 31
$wsMC.Range("A31").Value = "Wich of the follwing is a beeverage?"
$wsMC.Range("B31").Value = "Water"
$wsMC.Range("C31").Value = "Juicc"
$wsMC.Range("D31").Value = "Bread"
$wsMC.Range("E31").Value = "Cake"
$wsMC.Range("F31").Value = "Rice"
$wsMC.Range("G31").Value = "b"
$wsMC.Range("H31").Value = "drink1.png"
$wsMC.Range("I31").Value = "food"

# Row 32
$wsMC.Range("A32").Value = "Which animal can fly?"
$wsMC.Range("B32").Value = "Snake"
$wsMC.Range("C32").Value = "Dog"
$wsMC.Range("D32").Value = "Bird"
$wsMC.Range("E32").Value = "Cat"
$wsMC.Range("F32").Value = "Elephant"
$wsMC.Range("G32").Value = "c"
$wsMC.Range("H32").Value = "bird3.png"
$wsMC.Range("I32").Value = "animal"

# Row 33
$wsMC.Range("A33").Value = "Which animal can flyy?"
$wsMC.Range("B33").Value = "Fish"
$wsMC.Range("C33").Value = "Frog"
$wsMC.Range("D33").Value = "Bird"
$wsMC.Range("E33").Value = "Camel"
$wsMC.Range("F33").Value = "None"
$wsMC.Range("G33").Value = "c"
$wsMC.Range("H33").Value = "bird4.png"
$wsMC.Range("I33").Value = "animal"

# Row 34
$wsMC.Range("A34").Value = "How many corner does a square has?"
$wsMC.Range("B34").Value = 1
$wsMC.Range("C34").Value = 2
$wsMC.Range("D34").Value = 3
$wsMC.Range("E34").Value = 4
$wsMC.Range("F34").Value = 6
$wsMC.Range("G34").Value = "d"
$wsMC.Range("H34").Value = "square4.png"
$wsMC.Range("I34").Value = "maths"

# Row 35
$wsMC.Range("A35").Value = "How many cornes in square?"
$wsMC.Range("B35").Value = "One"
$wsMC.Range("C35").Value = "Two"
$wsMC.Range("D35").Value = "Three"
$wsMC.Range("E35").Value = "Four"
$wsMC.Range("F35").Value = "Five"
$wsMC.Range("G35").Value = "d"
$wsMC.Range("H35").Value = "square5.png"
$wsMC.Range("I35").Value = "maths"

# Row 36
$wsMC.Range("A36").Value = "Where is Chaoo Phara River?"
$wsMC.Range("B36").Value = "Bkk"
$wsMC.Range("C36").Value = "Russia"
$wsMC.Range("D36").Value = "England"
$wsMC.Range("E36").Value = "Your Home"
$wsMC.Range("F36").Value = "Germany"
$wsMC.Range("G36").Value = "a"
$wsMC.Range("H36").Value = "river2.png"
$wsMC.Range("I36").Value = "places"

# --- style: wrap text + vertical-center for all new MultipleChoice cells ---
$rngMC = $wsMC.Range("A19:I36")
$rngMC.WrapText = $true
$rngMC.VerticalAlignment = -4108

# --- row heights for rows whose wrapped text spans two lines ---
$wsMC.Rows.Item(19).RowHeight = 28.8
$wsMC.Rows.Item(25).RowHeight = 28.8
$wsMC.Rows.Item(26).RowHeight = 28.8
$wsMC.Rows.Item(34).RowHeight = 28.8
$wsMC.Rows.Item(35).RowHeight = 28.8

# --- widen column H (image filename) to fit the new, longer names ---
$wsMC.Columns.Item(8).ColumnWidth = 11.0533854166667

# --- Matching: append pair rows 10-33 ---
$wsM = $wb.Worksheets.Item("Matching")

# Row 10
$wsM.Range("A10").Value = "Something soft."
$wsM.Range("B10").Value = "Bread"
$wsM.Range("C10").Value = "food"
$wsM.Range("D10").Value = "5.png"

# Row 11
$wsM.Range("A11").Value = "Something crunchy."
$wsM.Range("B11").Value = "Chips"
$wsM.Range("C11").Value = "food"
$wsM.Range("D11").Value = "5.png"

# Row 12
$wsM.Range("A12").Value = "Something hot."
$wsM.Range("B12").Value = "Fire"
$wsM.Range("C12").Value = "temperature"
$wsM.Range("D12").Value = "3.png"

# Row 13
$wsM.Range("A13").Value = "Something cold."
$wsM.Range("B13").Value = "Ice"
$wsM.Range("C13").Value = "temperature"
$wsM.Range("D13").Value = "3.png"

# Row 14
$wsM.Range("A14").Value = "Something green."
$wsM.Range("B14").Value = "Spinach"
$wsM.Range("C14").Value = "vegetables"
$wsM.Range("D14").Value = "6.png"

# Row 15
$wsM.Range("A15").Value = "Something red."
$wsM.Range("B15").Value = "Tomato"
$wsM.Range("C15").Value = "vegetables"
$wsM.Range("D15").Value = "6.png"

# Row 16
$wsM.Range("A16").Value = "Something yellow."
$wsM.Range("B16").Value = "Banana"
$wsM.Range("C16").Value = "fruit"
$wsM.Range("D16").Value = "7.png"

# Row 17
$wsM.Range("A17").Value = "Something purple."
$wsM.Range("B17").Value = "Grapes"
$wsM.Range("C17").Value = "fruit"
$wsM.Range("D17").Value = "7.png"

# Row 18
$wsM.Range("A18").Value = "Something round."
$wsM.Range("B18").Value = "Ball"
$wsM.Range("C18").Value = "object"
$wsM.Range("D18").Value = "8.png"

# Row 19
$wsM.Range("A19").Value = "Something square."
$wsM.Range("B19").Value = "Dice"
$wsM.Range("C19").Value = "object"
$wsM.Range("D19").Value = "8.png"

# Row 20
$wsM.Range("A20").Value = "Something that fly."
$wsM.Range("B20").Value = "Bird"
$wsM.Range("C20").Value = "animal"
$wsM.Range("D20").Value = "9.png"

# Row 21
$wsM.Range("A21").Value = "Something that swim."
$wsM.Range("B21").Value = "Fish"
$wsM.Range("C21").Value = "animal"
$wsM.Range("D21").Value = "9.png"

# Row 22
$wsM.Range("A22").Value = "Something that bark."
$wsM.Range("B22").Value = "Dog"
$wsM.Range("C22").Value = "animal"
$wsM.Range("D22").Value = "10.png"

# Row 23
$wsM.Range("A23").Value = "Something that meow."
$wsM.Range("B23").Value = "Cat"
$wsM.Range("C23").Value = "animal"
$wsM.Range("D23").Value = "10.png"

# Row 24
$wsM.Range("A24").Value = "Something used to write."
$wsM.Range("B24").Value = "Pen"
$wsM.Range("C24").Value = "object"
$wsM.Range("D24").Value = "11.png"

# Row 25
$wsM.Range("A25").Value = "Something used to cut."
$wsM.Range("B25").Value = "Scissors"
$wsM.Range("C25").Value = "object"
$wsM.Range("D25").Value = "11.png"

# Row 26
$wsM.Range("A26").Value = "Something used to eat soup."
$wsM.Range("B26").Value = "Spoon"
$wsM.Range("C26").Value = "utensil"
$wsM.Range("D26").Value = "12.png"

# Row 27
$wsM.Range("A27").Value = "Something used to eat rice."
$wsM.Range("B27").Value = "Fork"
$wsM.Range("C27").Value = "utensil"
$wsM.Range("D27").Value = "12.png"

# Row 28
$wsM.Range("A28").Value = "Something you wear on foot."
$wsM.Range("B28").Value = "Shoes"
$wsM.Range("C28").Value = "clothing"
$wsM.Range("D28").Value = "13.png"

# Row 29
$wsM.Range("A29").Value = "Something you wear on head."
$wsM.Range("B29").Value = "Hat"
$wsM.Range("C29").Value = "clothing"
$wsM.Range("D29").Value = "13.png"

# Row 30
$wsM.Range("A30").Value = "Something that shine at night."
$wsM.Range("B30").Value = "Moon"
$wsM.Range("C30").Value = "nature"
$wsM.Range("D30").Value = "14.png"

# Row 31
$wsM.Range("A31").Value = "Something that shine at day."
$wsM.Range("B31").Value = "Sun"
$wsM.Range("C31").Value = "nature"
$wsM.Range("D31").Value = "14.png"

# Row 32
$wsM.Range("A32").Value = "Something we drink in morning."
$wsM.Range("B32").Value = "Coffee"
$wsM.Range("C32").Value = "drink"
$wsM.Range("D32").Value = "15.png"

# Row 33
$wsM.Range("A33").Value = "Something we drink when hot."
$wsM.Range("B33").Value = "Water"
$wsM.Range("C33").Value = "drink"
$wsM.Range("D33").Value = "15.png"

# --- style: wrap text + vertical-center for all new Matching cells ---
$rngM = $wsM.Range("A10:D33")
$rngM.WrapText = $true
$rngM.VerticalAlignment = -4108

# --- final selections / active sheet (matches the saved view state) ---
$wsInfo = $wb.Worksheets.Item("Info")
$wsTF = $wb.Worksheets.Item("TrueFalse")
$wsWQ = $wb.Worksheets.Item("WrittenQuestion")

$wsTF.Activate()
$wsTF.Range("G38").Select()

$wsWQ.Activate()
$wsWQ.Range("K7").Select()

$wsMC.Activate()
$wsMC.Range("E30").Select()

$wsM.Activate()
$wsM.Range("F30").Select()

